$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('empty')
$ws.Columns.Item(1).ColumnWidth = 0.0
$ws.Columns.Item(2).ColumnWidth = 0.25
$ws.Columns.Item(3).ColumnWidth = 0.5
$ws.Columns.Item(4).ColumnWidth = 0.75
$ws.Columns.Item(5).ColumnWidth = 1.0
$ws.Columns.Item(6).ColumnWidth = 1.25
$ws.Columns.Item(7).ColumnWidth = 1.5
$ws.Columns.Item(8).ColumnWidth = 1.75
$ws.Columns.Item(9).ColumnWidth = 2.0
$ws.Columns.Item(10).ColumnWidth = 2.25
$ws.Columns.Item(11).ColumnWidth = 2.5
$ws.Columns.Item(12).ColumnWidth = 2.75
$ws.Columns.Item(13).ColumnWidth = 3.0
$ws.Columns.Item(14).ColumnWidth = 3.25
$ws.Columns.Item(15).ColumnWidth = 3.5
$ws.Columns.Item(16).ColumnWidth = 3.75
$ws.Columns.Item(17).ColumnWidth = 4.0
$ws.Columns.Item(18).ColumnWidth = 4.25
$ws.Columns.Item(19).ColumnWidth = 4.5
$ws.Columns.Item(20).ColumnWidth = 4.75
$ws.Columns.Item(21).ColumnWidth = 5.0
$ws.Columns.Item(22).ColumnWidth = 5.25
$ws.Columns.Item(23).ColumnWidth = 5.5
$ws.Columns.Item(24).ColumnWidth = 5.75
$ws.Columns.Item(25).ColumnWidth = 6.0
$ws.Columns.Item(26).ColumnWidth = 6.25
$ws.Columns.Item(27).ColumnWidth = 6.5
$ws.Columns.Item(28).ColumnWidth = 6.75
$ws.Columns.Item(29).ColumnWidth = 7.0
$ws.Columns.Item(30).ColumnWidth = 7.25
$ws.Columns.Item(31).ColumnWidth = 7.5
$ws.Columns.Item(32).ColumnWidth = 7.75
$ws.Columns.Item(33).ColumnWidth = 8.0
$ws.Columns.Item(34).ColumnWidth = 8.25
$ws.Columns.Item(35).ColumnWidth = 8.5
$ws.Columns.Item(36).ColumnWidth = 8.75
$ws.Columns.Item(37).ColumnWidth = 9.0
$ws.Columns.Item(38).ColumnWidth = 9.25
$ws.Columns.Item(39).ColumnWidth = 9.5
$ws.Columns.Item(40).ColumnWidth = 9.75
$ws.Columns.Item(41).ColumnWidth = 10.0
$ws.Columns.Item(42).ColumnWidth = 10.25
$ws.Columns.Item(43).ColumnWidth = 10.5
$ws.Columns.Item(44).ColumnWidth = 10.75
$ws.Columns.Item(45).ColumnWidth = 11.0
$ws.Columns.Item(46).ColumnWidth = 11.25
$ws.Columns.Item(47).ColumnWidth = 11.5
$ws.Columns.Item(48).ColumnWidth = 11.75
$ws.Columns.Item(49).ColumnWidth = 12.0
$ws.Columns.Item(50).ColumnWidth = 12.25
$ws.Columns.Item(51).ColumnWidth = 12.5
$ws.Columns.Item(52).ColumnWidth = 12.75
$ws.Columns.Item(53).ColumnWidth = 13.0
$ws.Columns.Item(54).ColumnWidth = 13.25
$ws.Columns.Item(55).ColumnWidth = 13.5
$ws.Columns.Item(56).ColumnWidth = 13.75
$ws.Columns.Item(57).ColumnWidth = 14.0
$ws.Columns.Item(58).ColumnWidth = 14.25
$ws.Columns.Item(59).ColumnWidth = 14.5
$ws.Columns.Item(60).ColumnWidth = 14.75
$ws.Columns.Item(61).ColumnWidth = 15.0
$ws.Columns.Item(62).ColumnWidth = 15.25
$ws.Columns.Item(63).ColumnWidth = 15.5
$ws.Columns.Item(64).ColumnWidth = 15.75
$ws.Columns.Item(65).ColumnWidth = 16.0
$ws.Columns.Item(66).ColumnWidth = 16.25
$ws.Columns.Item(67).ColumnWidth = 16.5
$ws.Columns.Item(68).ColumnWidth = 16.75
$ws.Columns.Item(69).ColumnWidth = 17.0
$ws.Columns.Item(70).ColumnWidth = 17.25
$ws.Columns.Item(71).ColumnWidth = 17.5
$ws.Columns.Item(72).ColumnWidth = 17.75
$ws.Columns.Item(73).ColumnWidth = 18.0
$ws.Columns.Item(74).ColumnWidth = 18.25
$ws.Columns.Item(75).ColumnWidth = 18.5
$ws.Columns.Item(76).ColumnWidth = 18.75
$ws.Columns.Item(77).ColumnWidth = 19.0
$ws.Columns.Item(78).ColumnWidth = 19.25
$ws.Columns.Item(79).ColumnWidth = 19.5
$ws.Columns.Item(80).ColumnWidth = 19.75